$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Data Transaksi")

# Translate header cells from Indonesian to English
$ws.Range("I1").Value = "Created At"
$ws.Range("J1").Value = "Created By"
$ws.Range("C1").Value = "Branch"

# Update the active selection/view to reflect the saved state
$ws.Range("C2").Select()
